$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell style cloned from a size-11 "Normal" (as happens when pasting
# from another workbook whose default Normal style differs from this one).
$newStyle = $wb.Styles.Add("Normal 2")
$newStyle.Font.Size = 11

# New "comments" column with the article citation notes.
$ws.Range("G1").Value = "comments"
$ws.Range("G2").Value = "These data are based on Figure 5 in the article:"
$ws.Range("G3").Value = "The influence of soccer shoe design on playing performance: a series of"
$ws.Range("G4").Value = "biomechanical studies"
$ws.Range("G5").Value = "Ewald M. Hennig and Thorsten Sterzing"
$ws.Range("G6").Value = ""
$ws.Range("G7").Value = "Hennig, Ewald M. and Sterzing, Thorsten(2010) 'The influence of soccer shoe design on playing"
$ws.Range("G8").Value = "performance: a series of biomechanical studies', Footwear Science, 2: 1, 3 — 11"

$ws.Range("G1:G8").Style = "Normal 2"

$ws.Range("G1:G8").Select()
